# Sync cryptos.xlsx price/volume snapshot to the latest scrape, and fix
# up the coin ordering for rows 15-23 (CoinExToken moved up into the
# block; the rest of that block shifted down one row) - GitHub Actions
# refresh job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data cell on the sheet is stored as text (coin names, coin-page
# URLs, and numeric-looking Price/Volume strings such as "308.17" or
# "-3.99%"). `IsText` marks the numeric-looking ones so we can pre-format
# them as Text ("@") before writing - otherwise Excel would silently
# coerce the string into a Number/Percentage value.
$updates = @(
    @{ Cell = "D2"; Value = "308.17"; IsText = $true },
    @{ Cell = "E2"; Value = "-3.99%"; IsText = $true },
    @{ Cell = "D3"; Value = "39.59"; IsText = $true },
    @{ Cell = "E3"; Value = "-7.45%"; IsText = $true },
    @{ Cell = "D4"; Value = "5.103"; IsText = $true },
    @{ Cell = "E4"; Value = "-1.94%"; IsText = $true },
    @{ Cell = "D5"; Value = "0.07692"; IsText = $true },
    @{ Cell = "E5"; Value = "-5.85%"; IsText = $true },
    @{ Cell = "E6"; Value = "-1.85%"; IsText = $true },
    @{ Cell = "D7"; Value = "1.607"; IsText = $true },
    @{ Cell = "E7"; Value = "-11.24%"; IsText = $true },
    @{ Cell = "D8"; Value = "0.8968"; IsText = $true },
    @{ Cell = "E8"; Value = "-3.96%"; IsText = $true },
    @{ Cell = "E9"; Value = "-9.61%"; IsText = $true },
    @{ Cell = "D10"; Value = "0.1733"; IsText = $true },
    @{ Cell = "E10"; Value = "-6.45%"; IsText = $true },
    @{ Cell = "D11"; Value = "0.09027"; IsText = $true },
    @{ Cell = "E11"; Value = "-3.59%"; IsText = $true },
    @{ Cell = "D12"; Value = "0.04422"; IsText = $true },
    @{ Cell = "E12"; Value = "-6.18%"; IsText = $true },
    @{ Cell = "D13"; Value = "0.1053"; IsText = $true },
    @{ Cell = "E13"; Value = "-0.41%"; IsText = $true },
    @{ Cell = "D14"; Value = "0.001257"; IsText = $true },
    @{ Cell = "E14"; Value = "-2.82%"; IsText = $true },
    @{ Cell = "B15"; Value = "CoinExToken"; IsText = $false },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; IsText = $false },
    @{ Cell = "D15"; Value = "0.04146"; IsText = $true },
    @{ Cell = "E15"; Value = "-0.07%"; IsText = $true },
    @{ Cell = "B16"; Value = "TigerCash"; IsText = $false },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; IsText = $false },
    @{ Cell = "D16"; Value = "0.005921"; IsText = $true },
    @{ Cell = "E16"; Value = "-0.33%"; IsText = $true },
    @{ Cell = "B17"; Value = "UpBots"; IsText = $false },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"; IsText = $false },
    @{ Cell = "D17"; Value = "0.007491"; IsText = $true },
    @{ Cell = "E17"; Value = "2,411.68%"; IsText = $true },
    @{ Cell = "B18"; Value = "LEO"; IsText = $false },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; IsText = $false },
    @{ Cell = "D18"; Value = "3.353"; IsText = $true },
    @{ Cell = "E18"; Value = "-0.07%"; IsText = $true },
    @{ Cell = "B19"; Value = "BTSEToken"; IsText = $false },
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; IsText = $false },
    @{ Cell = "D19"; Value = "2.419"; IsText = $true },
    @{ Cell = "E19"; Value = "-3.53%"; IsText = $true },
    @{ Cell = "B20"; Value = "BitpandaEcosystemToken"; IsText = $false },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; IsText = $false },
    @{ Cell = "D20"; Value = "0.3317"; IsText = $true },
    @{ Cell = "E20"; Value = "-0.93%"; IsText = $true },
    @{ Cell = "B21"; Value = "MCDex"; IsText = $false },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"; IsText = $false },
    @{ Cell = "D21"; Value = "7.066"; IsText = $true },
    @{ Cell = "E21"; Value = "-5.50%"; IsText = $true },
    @{ Cell = "B22"; Value = "ProBitToken"; IsText = $false },
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"; IsText = $false },
    @{ Cell = "D22"; Value = "0.1349"; IsText = $true },
    @{ Cell = "E22"; Value = "-2.31%"; IsText = $true },
    @{ Cell = "B23"; Value = "ZBToken"; IsText = $false },
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"; IsText = $false },
    @{ Cell = "D23"; Value = "0.2763"; IsText = $true },
    @{ Cell = "E23"; Value = "8.43%"; IsText = $true },
    @{ Cell = "D24"; Value = "0.001209"; IsText = $true },
    @{ Cell = "E24"; Value = "-2.98%"; IsText = $true },
    @{ Cell = "D25"; Value = "0.004059"; IsText = $true },
    @{ Cell = "E25"; Value = "-5.60%"; IsText = $true },
    @{ Cell = "D26"; Value = "0.0001303"; IsText = $true },
    @{ Cell = "E26"; Value = "8.47%"; IsText = $true },
    @{ Cell = "D38"; Value = "0.02342"; IsText = $true },
    @{ Cell = "E38"; Value = "-13.27%"; IsText = $true },
    @{ Cell = "D39"; Value = "0.05182"; IsText = $true },
    @{ Cell = "E39"; Value = "-6.50%"; IsText = $true },
    @{ Cell = "D40"; Value = "0.007923"; IsText = $true },
    @{ Cell = "E40"; Value = "-1.52%"; IsText = $true },
    @{ Cell = "D41"; Value = "0.1318"; IsText = $true },
    @{ Cell = "E41"; Value = "-5.60%"; IsText = $true },
    @{ Cell = "D42"; Value = "0.006296"; IsText = $true },
    @{ Cell = "E42"; Value = "-3.84%"; IsText = $true },
    @{ Cell = "D43"; Value = "0.001954"; IsText = $true },
    @{ Cell = "E43"; Value = "-6.41%"; IsText = $true },
    @{ Cell = "D44"; Value = "0.008232"; IsText = $true },
    @{ Cell = "E44"; Value = "-0.36%"; IsText = $true },
    @{ Cell = "D45"; Value = "0.3333"; IsText = $true },
    @{ Cell = "E45"; Value = "-4.60%"; IsText = $true },
    @{ Cell = "E46"; Value = "-5.87%"; IsText = $true },
    @{ Cell = "D47"; Value = "0.00000000752"; IsText = $true },
    @{ Cell = "E47"; Value = "0.12%"; IsText = $true },
    @{ Cell = "E48"; Value = "98.23%"; IsText = $true },
    @{ Cell = "D49"; Value = "0.003501"; IsText = $true },
    @{ Cell = "E49"; Value = "4.83%"; IsText = $true },
    @{ Cell = "E50"; Value = "0.12%"; IsText = $true },
    @{ Cell = "E51"; Value = "0.12%"; IsText = $true }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.IsText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
